# Bond dates update: recompute "Dni od poprzedniej wypłaty" (col G) and
# "Dni do następnej wypłaty" (col I) as of the new reference date
# 2023-09-22 (Excel serial 45191), one day later than the previous
# reference date (2023-09-21 / serial 45190) baked into the sheet.
#
#   G (days since previous payout) = RefDate - F (previous payout date)
#   I (days until next payout)     = H (next payout date) - RefDate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$refDate = 45191

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2   # column F
    if ($fVal -ne $null) {
        $ws.Cells.Item($r, 7).Value2 = $refDate - $fVal   # column G
    }

    $hVal = $ws.Cells.Item($r, 8).Value2   # column H
    if ($hVal -ne $null) {
        $ws.Cells.Item($r, 9).Value2 = $hVal - $refDate   # column I
    }
}
